$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "64.242.52"
Set-TextValue $ws.Range("E2") "  +0.10%  "
Set-TextValue $ws.Range("D3") "3.493.04"
Set-TextValue $ws.Range("E3") "  -0.83%  "
Set-TextValue $ws.Range("E4") "  +0.06%  "
Set-TextValue $ws.Range("D5") "587.19"
Set-TextValue $ws.Range("E5") "  +0.11%  "
Set-TextValue $ws.Range("D6") "134.13"
Set-TextValue $ws.Range("E6") "  +0.54%  "
Set-TextValue $ws.Range("E7") "  -0.02%  "
Set-TextValue $ws.Range("E8") "  +0.09%  "
Set-TextValue $ws.Range("E9") "  -0.40%  "
Set-TextValue $ws.Range("E10") "  +1.83%  "
Set-TextValue $ws.Range("E11") "  +1.77%  "
Set-TextValue $ws.Range("D12") "4.086.33"
Set-TextValue $ws.Range("E12") "  -0.94%  "
Set-TextValue $ws.Range("E13") "  +1.03%  "
Set-TextValue $ws.Range("E14") "  +1.20%  "
Set-TextValue $ws.Range("D15") "3.490.42"
Set-TextValue $ws.Range("E15") "  -1.19%  "
Set-TextValue $ws.Range("B16") "WrappedBTC"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D16") "64.326.13"
Set-TextValue $ws.Range("E16") "  +0.22%  "
Set-TextValue $ws.Range("B17") "Avalanche"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue $ws.Range("D17") "25.70"
Set-TextValue $ws.Range("E17") "  -6.98%  "
Set-TextValue $ws.Range("E18") "  +0.90%  "
Set-TextValue $ws.Range("E19") "  +2.12%  "
Set-TextValue $ws.Range("D20") "13.58"
Set-TextValue $ws.Range("E20") "  -3.25%  "
Set-TextValue $ws.Range("D21") "393.96"
Set-TextValue $ws.Range("E21") "  +2.19%  "
Set-TextValue $ws.Range("E22") "  -0.84%  "
Set-TextValue $ws.Range("D23") "3.631.88"
Set-TextValue $ws.Range("E23") "  -0.90%  "
Set-TextValue $ws.Range("E24") "  +0.99%  "
Set-TextValue $ws.Range("E25") "  +0.04%  "
Set-TextValue $ws.Range("E26") "  +0.31%  "
Set-TextValue $ws.Range("E27") "  +0.14%  "
Set-TextValue $ws.Range("D28") "0.994"
Set-TextValue $ws.Range("E28") "  -0.58%  "
Set-TextValue $ws.Range("E29") "  -2.15%  "
Set-TextValue $ws.Range("E30") "  +0.09%  "
Set-TextValue $ws.Range("E31") "  -2.60%  "
Set-TextValue $ws.Range("D32") "1.47"
Set-TextValue $ws.Range("E32") "  -5.90%  "
Set-TextValue $ws.Range("D33") "3.514.13"
Set-TextValue $ws.Range("E33") "  -0.53%  "
Set-TextValue $ws.Range("E34") "  +3.66%  "
Set-TextValue $ws.Range("E35") "  +0.02%  "
Set-TextValue $ws.Range("E36") "  -0.87%  "
Set-TextValue $ws.Range("E37") "  -4.60%  "
Set-TextValue $ws.Range("E38") "  -0.68%  "
Set-TextValue $ws.Range("D39") "6.87"
Set-TextValue $ws.Range("E39") "  -0.83%  "
Set-TextValue $ws.Range("D40") "166.45"
Set-TextValue $ws.Range("E40") "  +3.38%  "
Set-TextValue $ws.Range("D41") "0.0777"
Set-TextValue $ws.Range("E41") "  -1.46%  "
Set-TextValue $ws.Range("E42") "  -1.07%  "
Set-TextValue $ws.Range("E43") "  +0.10%  "
Set-TextValue $ws.Range("D44") "25.28"
Set-TextValue $ws.Range("E44") "  -4.68%  "
Set-TextValue $ws.Range("E45") "  -0.93%  "
Set-TextValue $ws.Range("E46") "  +2.02%  "
Set-TextValue $ws.Range("E47") "  -4.04%  "
Set-TextValue $ws.Range("D48") "2.460.62"
Set-TextValue $ws.Range("E48") "  -0.63%  "
Set-TextValue $ws.Range("D49") "6.74"
Set-TextValue $ws.Range("E49") "  -0.98%  "
Set-TextValue $ws.Range("E50") "  -1.59%  "
Set-TextValue $ws.Range("E51") "  -1.23%  "
